# "MP Taxonomy with River Data"
# Adds a new Color taxonomy entry (Gray / Grey alias, sourced from a new
# reference) plus two additional color rows (Purple, Unknown) to Sheet1.
#
# Cell-entry order below matches the order new values were typed into the
# sheet (and therefore the order they land in the shared-string table):
# Purple, Unknown, Grey, Gray, then the new reference DOI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New color rows appended first (row 14 = Purple, row 15 = Unknown)
$ws.Range("C14").Value = "Purple"
$ws.Range("C15").Value = "Unknown"

# Then fill in row 13: a Gray/Grey entry with its supporting reference
$ws.Range("D13").Value = "Grey"
$ws.Range("C13").Value = "Gray"
$ws.Range("A13").Value = "doi.org/10.1016/j.envpol.2016.01.018"

# Leave the active selection on C14, matching the saved view state
[void]$ws.Range("C14").Select()
